$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "iaest-measure:formacion-jefe-explotacion"
$ws.Range("D2").Value = "sdmx-dimension:refArea"

$ws.Range("A3").Value = "medida"
$ws.Range("D3").Value = "dim"

$ws.Range("A4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-Municipio"

$ws.Range("A5").Clear()
